$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new detail rows -------------------------------------
# Old layout: detail rows 16-22 (normal style) + row 23 (thicker "last row"
# style). New layout adds 6 more detail rows, so the "last row" style needs
# to end up on row 29 and rows 23-28 need the normal detail-row formatting.
$ws.Range("B23:J28").Insert()

# Copy the normal detail-row formatting (from row 22) down into the 6 new
# rows (23-28) so they match rows 16-22's look (borders/number formats).
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J28").PasteSpecial(-4122)

# --- Header values -----------------------------------------------------------
$ws.Range("E11").Value = 436810
$ws.Range("C13").Value = 2

# --- Detail rows: new employee "ANGELICA MARIA RODRIGUEZ TORRES" -----------
$angelica = "45551082"
$angelicaName = "ANGELICA MARIA RODRIGUEZ TORRES"
$salario = 828116

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = $angelica
$ws.Range("D16").Value = $angelicaName
$ws.Range("E16").Value = "2004"
$ws.Range("F16").Value = 16562
$ws.Range("G16").Value = $salario

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = $angelica
$ws.Range("D17").Value = $angelicaName
$ws.Range("E17").Value = "2003"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = $salario

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = $angelica
$ws.Range("D18").Value = $angelicaName
$ws.Range("E18").Value = "2002"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = $salario

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = $angelica
$ws.Range("D19").Value = $angelicaName
$ws.Range("E19").Value = "2001"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = $salario

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = $angelica
$ws.Range("D20").Value = $angelicaName
$ws.Range("E20").Value = "1910"
$ws.Range("F20").Value = 33125
$ws.Range("G20").Value = $salario

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = $angelica
$ws.Range("D21").Value = $angelicaName
$ws.Range("E21").Value = "1909"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = $salario

# --- Detail rows: existing employee "CARMEN ROCIO BALLESTEROS FLOREZ" ------
# now with 8 periods of arrears instead of 1.
$carmen = "32905679"
$carmenName = "CARMEN ROCIO BALLESTEROS FLOREZ"

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = $carmen
$ws.Range("D22").Value = $carmenName
$ws.Range("E22").Value = "2006"
$ws.Range("F22").Value = 26500
$ws.Range("G22").Value = $salario

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = $carmen
$ws.Range("D23").Value = $carmenName
$ws.Range("E23").Value = "2005"
$ws.Range("F23").Value = 33125
$ws.Range("G23").Value = $salario

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = $carmen
$ws.Range("D24").Value = $carmenName
$ws.Range("E24").Value = "2004"
$ws.Range("F24").Value = 33125
$ws.Range("G24").Value = $salario

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = $carmen
$ws.Range("D25").Value = $carmenName
$ws.Range("E25").Value = "2003"
$ws.Range("F25").Value = 33125
$ws.Range("G25").Value = $salario

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = $carmen
$ws.Range("D26").Value = $carmenName
$ws.Range("E26").Value = "2002"
$ws.Range("F26").Value = 33125
$ws.Range("G26").Value = $salario

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = $carmen
$ws.Range("D27").Value = $carmenName
$ws.Range("E27").Value = "2001"
$ws.Range("F27").Value = 33125
$ws.Range("G27").Value = $salario

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = $carmen
$ws.Range("D28").Value = $carmenName
$ws.Range("E28").Value = "1910"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = $salario

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = $carmen
$ws.Range("D29").Value = $carmenName
$ws.Range("E29").Value = "1909"
$ws.Range("F29").Value = 31249
$ws.Range("G29").Value = $salario

# --- Column widths (Excel's bestFit autofit grew these once the new, wider
# values were entered) ------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.709635416666668
$ws.Columns.Item(3).ColumnWidth = 15.893229166666666
$ws.Columns.Item(5).ColumnWidth = 12.709635416666666
$ws.Columns.Item(6).ColumnWidth = 9.346354166666666
$ws.Columns.Item(7).ColumnWidth = 13.529947916666666
$ws.Columns.Item(8).ColumnWidth = 18.529947916666668
$ws.Columns.Item(9).ColumnWidth = 17.256510416666668
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
